# Updates cryptos list values per upstream data refresh (Fri Oct 18 16:34:44 UTC 2024, GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every target cell to Text format before writing so numeric-looking
# strings (e.g. "599.80", "0.350") are not auto-coerced into numbers by Excel,
# then restore the default "Normal" style so no spurious formatting diff remains.
$cellUpdates = @{
    'D2' = '68.615.29'
    'E2' = '  +2.00%  '
    'D3' = '2.641.22'
    'E3' = '  +0.98%  '
    'E4' = '  +0.02%  '
    'D5' = '599.80'
    'E5' = '  +1.37%  '
    'D6' = '154.48'
    'E6' = '  +2.07%  '
    'E7' = '  -0.01%  '
    'D8' = '0.545'
    'E8' = '  -1.27%  '
    'D9' = '2.640.47'
    'E9' = '  +0.95%  '
    'D10' = '0.138'
    'E10' = '  +13.09%  '
    'E11' = '  -0.51%  '
    'D12' = '5.25'
    'E12' = '  +0.96%  '
    'D13' = '0.350'
    'E13' = '  +1.03%  '
    'E14' = '  +0.52%  '
    'D15' = '0.0000189'
    'E15' = '  +5.72%  '
    'D16' = '3.120.69'
    'E16' = '  +0.65%  '
    'D17' = '68.494.27'
    'D18' = '2.636.56'
    'E18' = '  +0.52%  '
    'D19' = '11.43'
    'E19' = '  +3.29%  '
    'D20' = '365.28'
    'E20' = '  -0.51%  '
    'D21' = '7.43'
    'E21' = '  +0.45%  '
    'E22' = '  -0.85%  '
    'D23' = '4.87'
    'E23' = '  +0.45%  '
    'D24' = '2.08'
    'E24' = '  +1.68%  '
    'D25' = '73.06'
    'E25' = '  +10.40%  '
    'E26' = '  +0.12%  '
    'D27' = '9.96'
    'E27' = '  -1.17%  '
    'D28' = '2.790.86'
    'E28' = '  +1.11%  '
    'D29' = '0.0000105'
    'E29' = '  +4.69%  '
    'E30' = '  +0.02%  '
    'D31' = '578.91'
    'E31' = '  -0.86%  '
    'B32' = 'Fetch.AI'
    'C32' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D32' = '1.42'
    'E32' = '  +2.84%  '
    'B33' = 'InternetComputer(DFINITY)'
    'C33' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D33' = '7.97'
    'E33' = '  +3.77%  '
    'D34' = '1.85'
    'E34' = '  +2.60%  '
    'E35' = '  +0.01%  '
    'D36' = '0.128'
    'E36' = '  +3.17%  '
    'D37' = '1.54'
    'E37' = '  +2.51%  '
    'D38' = '160.25'
    'E38' = '  +2.44%  '
    'D39' = '19.29'
    'E39' = '  +1.47%  '
    'D40' = '1.92'
    'E40' = '  +1.31%  '
    'D41' = '0.367'
    'E41' = '  +0.41%  '
    'D42' = '5.37'
    'E42' = '  +2.73%  '
    'D43' = '2.66'
    'E43' = '  +4.22%  '
    'E44' = '  +5.34%  '
    'D45' = '0.0₆0321'
    'E45' = '  +8.69%  '
    'E46' = '  +0.02%  '
    'E47' = '  -0.76%  '
    'D48' = '155.75'
    'E48' = '  +0.25%  '
    'D49' = '3.72'
    'E49' = '  -0.16%  '
    'D50' = '22.06'
    'E50' = '  +1.28%  '
    'E51' = '  +0.30%  '
}

foreach ($ref in $cellUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$ref]
    $cell.Style = "Normal"
}
